# -----------------------------------------------------------------------
# Refresh the scraped-schedule data for Línea 141 (workbook: horarios-141)
# across all three sheets: LP1912, LP1912-215, 6203-6173.
#
# The scraper re-ran (new "Última actualización" timestamp 08:49:35,
# replacing 08:32:32) and produced a larger, re-sorted dataset; this
# script rewrites the header metadata and the affected data rows to
# match the refreshed feed.
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # LP1912
$ws2 = $wb.Worksheets.Item(2)   # LP1912-215
$ws3 = $wb.Worksheets.Item(3)   # 6203-6173

# =========================================================================
# Sheet 1: LP1912
# =========================================================================

# --- Header metadata ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 08:49:35"
$ws1.Cells.Item(3,1).Value = "Total filas: 130"

# --- Rows 55-56 and 87-89: reordered by arrival time (values swapped) ---
# --- Rows 98-127: refreshed scrape data (existing rows overwritten)   ---
# --- Rows 128-135: brand-new rows appended by the refreshed scrape    ---
$ws1.Cells.Item(55,1).Value = "05:23:05"
$ws1.Cells.Item(55,2).Value = "07:16"
$ws1.Cells.Item(55,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(55,4).Value = 113
$ws1.Cells.Item(55,5).Value = "LP1912"
$ws1.Cells.Item(56,1).Value = "06:52:34"
$ws1.Cells.Item(56,2).Value = "07:16"
$ws1.Cells.Item(56,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(56,4).Value = 24
$ws1.Cells.Item(56,5).Value = "LP1912"
$ws1.Cells.Item(87,1).Value = "07:49:14"
$ws1.Cells.Item(87,2).Value = "08:23"
$ws1.Cells.Item(87,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(87,4).Value = 34
$ws1.Cells.Item(87,5).Value = "LP1912"
$ws1.Cells.Item(88,1).Value = "08:02:29"
$ws1.Cells.Item(88,2).Value = "08:23"
$ws1.Cells.Item(88,3).Value = "215B_EL PATO"
$ws1.Cells.Item(88,4).Value = 21
$ws1.Cells.Item(88,5).Value = "LP1912"
$ws1.Cells.Item(89,1).Value = "07:18:13"
$ws1.Cells.Item(89,2).Value = "08:23"
$ws1.Cells.Item(89,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(89,4).Value = 65
$ws1.Cells.Item(89,5).Value = "LP1912"
$ws1.Cells.Item(98,1).Value = "08:49:35"
$ws1.Cells.Item(98,2).Value = "08:52"
$ws1.Cells.Item(98,3).Value = "10_OLMOS"
$ws1.Cells.Item(98,4).Value = 3
$ws1.Cells.Item(98,5).Value = "LP1912"
$ws1.Cells.Item(99,1).Value = "08:32:32"
$ws1.Cells.Item(99,2).Value = "08:53"
$ws1.Cells.Item(99,3).Value = "10_OLMOS"
$ws1.Cells.Item(99,4).Value = 21
$ws1.Cells.Item(99,5).Value = "LP1912"
$ws1.Cells.Item(100,1).Value = "08:49:35"
$ws1.Cells.Item(100,2).Value = "08:54"
$ws1.Cells.Item(100,3).Value = "17_ROMERO"
$ws1.Cells.Item(100,4).Value = 5
$ws1.Cells.Item(100,5).Value = "LP1912"
$ws1.Cells.Item(101,1).Value = "08:49:35"
$ws1.Cells.Item(101,2).Value = "09:01"
$ws1.Cells.Item(101,3).Value = "215A_EL PATO"
$ws1.Cells.Item(101,4).Value = 12
$ws1.Cells.Item(101,5).Value = "LP1912"
$ws1.Cells.Item(102,1).Value = "07:49:14"
$ws1.Cells.Item(102,2).Value = "09:02"
$ws1.Cells.Item(102,3).Value = "215A_EL PATO"
$ws1.Cells.Item(102,4).Value = 73
$ws1.Cells.Item(102,5).Value = "LP1912"
$ws1.Cells.Item(103,1).Value = "08:49:35"
$ws1.Cells.Item(103,2).Value = "09:03"
$ws1.Cells.Item(103,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(103,4).Value = 14
$ws1.Cells.Item(103,5).Value = "LP1912"
$ws1.Cells.Item(104,1).Value = "08:32:32"
$ws1.Cells.Item(104,2).Value = "09:04"
$ws1.Cells.Item(104,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(104,4).Value = 32
$ws1.Cells.Item(104,5).Value = "LP1912"
$ws1.Cells.Item(105,1).Value = "07:49:14"
$ws1.Cells.Item(105,2).Value = "09:04"
$ws1.Cells.Item(105,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(105,4).Value = 75
$ws1.Cells.Item(105,5).Value = "LP1912"
$ws1.Cells.Item(106,1).Value = "08:49:35"
$ws1.Cells.Item(106,2).Value = "09:05"
$ws1.Cells.Item(106,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(106,4).Value = 16
$ws1.Cells.Item(106,5).Value = "LP1912"
$ws1.Cells.Item(107,1).Value = "08:02:29"
$ws1.Cells.Item(107,2).Value = "09:08"
$ws1.Cells.Item(107,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(107,4).Value = 66
$ws1.Cells.Item(107,5).Value = "LP1912"
$ws1.Cells.Item(108,1).Value = "08:32:32"
$ws1.Cells.Item(108,2).Value = "09:09"
$ws1.Cells.Item(108,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(108,4).Value = 37
$ws1.Cells.Item(108,5).Value = "LP1912"
$ws1.Cells.Item(109,1).Value = "08:49:35"
$ws1.Cells.Item(109,2).Value = "09:10"
$ws1.Cells.Item(109,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(109,4).Value = 21
$ws1.Cells.Item(109,5).Value = "LP1912"
$ws1.Cells.Item(110,1).Value = "08:49:35"
$ws1.Cells.Item(110,2).Value = "09:10"
$ws1.Cells.Item(110,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(110,4).Value = 21
$ws1.Cells.Item(110,5).Value = "LP1912"
$ws1.Cells.Item(111,1).Value = "07:49:14"
$ws1.Cells.Item(111,2).Value = "09:11"
$ws1.Cells.Item(111,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(111,4).Value = 82
$ws1.Cells.Item(111,5).Value = "LP1912"
$ws1.Cells.Item(112,1).Value = "08:49:35"
$ws1.Cells.Item(112,2).Value = "09:16"
$ws1.Cells.Item(112,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(112,4).Value = 27
$ws1.Cells.Item(112,5).Value = "LP1912"
$ws1.Cells.Item(113,1).Value = "07:49:14"
$ws1.Cells.Item(113,2).Value = "09:17"
$ws1.Cells.Item(113,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(113,4).Value = 88
$ws1.Cells.Item(113,5).Value = "LP1912"
$ws1.Cells.Item(114,1).Value = "08:49:35"
$ws1.Cells.Item(114,2).Value = "09:21"
$ws1.Cells.Item(114,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(114,4).Value = 32
$ws1.Cells.Item(114,5).Value = "LP1912"
$ws1.Cells.Item(115,1).Value = "08:02:29"
$ws1.Cells.Item(115,2).Value = "09:21"
$ws1.Cells.Item(115,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(115,4).Value = 79
$ws1.Cells.Item(115,5).Value = "LP1912"
$ws1.Cells.Item(116,1).Value = "08:49:35"
$ws1.Cells.Item(116,2).Value = "09:22"
$ws1.Cells.Item(116,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(116,4).Value = 33
$ws1.Cells.Item(116,5).Value = "LP1912"
$ws1.Cells.Item(117,1).Value = "08:32:32"
$ws1.Cells.Item(117,2).Value = "09:22"
$ws1.Cells.Item(117,3).Value = "17_ROMERO"
$ws1.Cells.Item(117,4).Value = 50
$ws1.Cells.Item(117,5).Value = "LP1912"
$ws1.Cells.Item(118,1).Value = "08:02:29"
$ws1.Cells.Item(118,2).Value = "09:23"
$ws1.Cells.Item(118,3).Value = "17_ROMERO"
$ws1.Cells.Item(118,4).Value = 81
$ws1.Cells.Item(118,5).Value = "LP1912"
$ws1.Cells.Item(119,1).Value = "08:49:35"
$ws1.Cells.Item(119,2).Value = "09:23"
$ws1.Cells.Item(119,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(119,4).Value = 34
$ws1.Cells.Item(119,5).Value = "LP1912"
$ws1.Cells.Item(120,1).Value = "07:49:14"
$ws1.Cells.Item(120,2).Value = "09:24"
$ws1.Cells.Item(120,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(120,4).Value = 95
$ws1.Cells.Item(120,5).Value = "LP1912"
$ws1.Cells.Item(121,1).Value = "08:49:35"
$ws1.Cells.Item(121,2).Value = "09:32"
$ws1.Cells.Item(121,3).Value = "15_ABASTO"
$ws1.Cells.Item(121,4).Value = 43
$ws1.Cells.Item(121,5).Value = "LP1912"
$ws1.Cells.Item(122,1).Value = "08:49:35"
$ws1.Cells.Item(122,2).Value = "09:33"
$ws1.Cells.Item(122,3).Value = "10_OLMOS"
$ws1.Cells.Item(122,4).Value = 44
$ws1.Cells.Item(122,5).Value = "LP1912"
$ws1.Cells.Item(123,1).Value = "08:49:35"
$ws1.Cells.Item(123,2).Value = "09:34"
$ws1.Cells.Item(123,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(123,4).Value = 45
$ws1.Cells.Item(123,5).Value = "LP1912"
$ws1.Cells.Item(124,1).Value = "08:49:35"
$ws1.Cells.Item(124,2).Value = "09:42"
$ws1.Cells.Item(124,3).Value = "215C_EL PATO"
$ws1.Cells.Item(124,4).Value = 53
$ws1.Cells.Item(124,5).Value = "LP1912"
$ws1.Cells.Item(125,1).Value = "08:49:35"
$ws1.Cells.Item(125,2).Value = "09:43"
$ws1.Cells.Item(125,3).Value = "14_ABASTO"
$ws1.Cells.Item(125,4).Value = 54
$ws1.Cells.Item(125,5).Value = "LP1912"
$ws1.Cells.Item(126,1).Value = "07:49:14"
$ws1.Cells.Item(126,2).Value = "09:44"
$ws1.Cells.Item(126,3).Value = "14_ABASTO"
$ws1.Cells.Item(126,4).Value = 115
$ws1.Cells.Item(126,5).Value = "LP1912"
$ws1.Cells.Item(127,1).Value = "08:32:32"
$ws1.Cells.Item(127,2).Value = "09:47"
$ws1.Cells.Item(127,3).Value = "10_OLMOS"
$ws1.Cells.Item(127,4).Value = 75
$ws1.Cells.Item(127,5).Value = "LP1912"
$ws1.Cells.Item(128,1).Value = "08:49:35"
$ws1.Cells.Item(128,2).Value = "09:52"
$ws1.Cells.Item(128,3).Value = "15_ABASTO"
$ws1.Cells.Item(128,4).Value = 63
$ws1.Cells.Item(128,5).Value = "LP1912"
$ws1.Cells.Item(129,1).Value = "08:49:35"
$ws1.Cells.Item(129,2).Value = "09:53"
$ws1.Cells.Item(129,3).Value = "10_OLMOS"
$ws1.Cells.Item(129,4).Value = 64
$ws1.Cells.Item(129,5).Value = "LP1912"
$ws1.Cells.Item(130,1).Value = "08:49:35"
$ws1.Cells.Item(130,2).Value = "10:10"
$ws1.Cells.Item(130,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(130,4).Value = 81
$ws1.Cells.Item(130,5).Value = "LP1912"
$ws1.Cells.Item(131,1).Value = "08:32:32"
$ws1.Cells.Item(131,2).Value = "10:12"
$ws1.Cells.Item(131,3).Value = "15_ABASTO"
$ws1.Cells.Item(131,4).Value = 100
$ws1.Cells.Item(131,5).Value = "LP1912"
$ws1.Cells.Item(132,1).Value = "08:49:35"
$ws1.Cells.Item(132,2).Value = "10:21"
$ws1.Cells.Item(132,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(132,4).Value = 92
$ws1.Cells.Item(132,5).Value = "LP1912"
$ws1.Cells.Item(133,1).Value = "08:49:35"
$ws1.Cells.Item(133,2).Value = "10:26"
$ws1.Cells.Item(133,3).Value = "215A_EL PATO"
$ws1.Cells.Item(133,4).Value = 97
$ws1.Cells.Item(133,5).Value = "LP1912"
$ws1.Cells.Item(134,1).Value = "08:49:35"
$ws1.Cells.Item(134,2).Value = "10:42"
$ws1.Cells.Item(134,3).Value = "17_ROMERO"
$ws1.Cells.Item(134,4).Value = 113
$ws1.Cells.Item(134,5).Value = "LP1912"
$ws1.Cells.Item(135,1).Value = "08:49:35"
$ws1.Cells.Item(135,2).Value = "10:43"
$ws1.Cells.Item(135,3).Value = "14_ABASTO"
$ws1.Cells.Item(135,4).Value = 114
$ws1.Cells.Item(135,5).Value = "LP1912"

# =========================================================================
# Sheet 2: LP1912-215
# =========================================================================

# --- Header metadata ---
$ws2.Cells.Item(2,1).Value = "Última actualización: 08:49:35"

# --- Rows 24, 26, 27: Hora_Scrap + Minutos refreshed ---
$ws2.Cells.Item(24,1).Value = "08:49:35"
$ws2.Cells.Item(24,4).Value = 12
$ws2.Cells.Item(26,1).Value = "08:49:35"
$ws2.Cells.Item(26,4).Value = 53
$ws2.Cells.Item(27,1).Value = "08:49:35"
$ws2.Cells.Item(27,4).Value = 97

# =========================================================================
# Sheet 3: 6203-6173
# =========================================================================

# --- Header metadata ---
$ws3.Cells.Item(2,1).Value = "Última actualización: 08:49:35"
$ws3.Cells.Item(3,1).Value = "Total filas: 20"

# --- New row inserted at position 23 (pushes old rows 23.. down by 1) ---
$ws3.Rows.Item(23).Insert()

$ws3.Cells.Item(23,1).Value = "08:49:35"
$ws3.Cells.Item(23,2).Value = "09:08"
$ws3.Cells.Item(23,3).Value = "215D_LA PLATA"
$ws3.Cells.Item(23,4).Value = 19
$ws3.Cells.Item(23,5).Value = "L6203"

# --- Old row 24 (now shifted to row 25): Hora_Scrap + Minutos refreshed ---
$ws3.Cells.Item(25,1).Value = "08:49:35"
$ws3.Cells.Item(25,4).Value = 74
